# Update the stratigraphic color assignments for the Wolfcamp XY/A/B/C/D
# rows (H14:H20) - the old teal/olive/orange swatches are replaced with a
# descending red-shade palette, and nudge the remembered cell selection
# from I30 to I31 (an incidental artifact of the editing session).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H14").Value = "#CC0000"
$ws.Range("H15").Value = "#B20000"
$ws.Range("H16").Value = "#990000"
$ws.Range("H17").Value = "#800000"
$ws.Range("H18").Value = "#660000"
$ws.Range("H19").Value = "#4D0000"
$ws.Range("H20").Value = "#330000"

$ws.Range("I31").Select()
